$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to target "ECs" now targets "FAPs" (intra-cluster signalling),
# and the row that used to target "FAPs" now targets "MuSCs". The numeric TPM-derived
# columns (K:T) were recomputed with the new TPM values; the last data row (old row 4,
# FAPs -> MuSCs duplicate) is removed entirely since the new table only has two rows.

# Remove the now-obsolete last data row (old row 4).
$ws.Rows(4).Delete()

# Row 2: target cluster ECs -> FAPs, and recomputed receptor/edge metrics.
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.559997666666667
$ws.Range("N2").Value = 7.679993000000001
$ws.Range("O2").Value = 0.3067691265831628
$ws.Range("P2").Value = 0.3067691265831629
$ws.Range("Q2").Value = 0.1259228718931111
$ws.Range("R2").Value = 1.133305847038
$ws.Range("S2").Value = 0.3067691265831628
$ws.Range("T2").Value = 0.3067691265831629

# Row 3: target cluster FAPs -> MuSCs, and recomputed receptor/edge metrics.
$ws.Range("D3").Value = "MuSCs"
$ws.Range("M3").Value = 5.785032666666666
$ws.Range("N3").Value = 17.355098
$ws.Range("O3").Value = 0.6932308734168371
$ws.Range("P3").Value = 0.6932308734168372
$ws.Range("Q3").Value = 0.2845580434964444
$ws.Range("R3").Value = 2.561022391468
$ws.Range("S3").Value = 0.6932308734168371
$ws.Range("T3").Value = 0.6932308734168372
